$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.926.11"
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").Value = "2.571.78"
$ws.Range("E3").Value = "  +2.73%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.54"
$ws.Range("E5").Value = "  +2.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.06"
$ws.Range("E6").Value = "  +4.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  +1.55%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.550"
$ws.Range("E9").Value = "  +1.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.48"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.65"
$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("E13").Value = "  +7.67%  "

$ws.Range("D14").Value = "2.572.83"
$ws.Range("E14").Value = "  +3.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.887"
$ws.Range("E15").Value = "  +3.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.39"
$ws.Range("E16").Value = "  +3.11%  "

$ws.Range("D17").Value = "42.878.93"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.04"
$ws.Range("E18").Value = "  +6.60%  "

$ws.Range("E19").Value = "  +3.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  +3.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.10"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.81"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.98"
$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "28.56"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("E27").Value = "  +3.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.22"
$ws.Range("E28").Value = "  +4.51%  "

$ws.Range("E29").Value = "  -3.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.65"
$ws.Range("E31").Value = "  +3.33%  "

$ws.Range("E32").Value = "  +0.69%  "

$ws.Range("E33").Value = "  +2.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0815"
$ws.Range("E34").Value = "  +3.06%  "

$ws.Range("E35").Value = "  -1.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.35"
$ws.Range("E36").Value = "  +12.22%  "

$ws.Range("E37").Value = "  +1.71%  "

$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.57"
$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.44"
$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.90"
$ws.Range("E41").Value = "  +2.72%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0311"
$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("E43").Value = "  +28.99%  "

$ws.Range("D44").Value = "2.068.35"
$ws.Range("E44").Value = "  +3.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("E46").Value = "  +6.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.38"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "77.07"
$ws.Range("E48").Value = "  +13.96%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.21"
$ws.Range("E49").Value = "  +3.95%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.821.98"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("E51").Value = "  +3.67%  "
